$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.472.31'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.904.39'
$ws.Range('E3').Value = '  +2.58%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.28'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.634'
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.98'
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('E9').Value = '  +2.74%  '
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '2.179.96'
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '12.32'
$ws.Range('E13').Value = '  +7.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.693'
$ws.Range('E14').Value = '  +1.78%  '
$ws.Range('D15').Value = '1.895.00'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('E16').Value = '  +3.42%  '
$ws.Range('D17').Value = '35.486.09'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '71.77'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').Value = '0.0₃0823'
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '243.19'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.60'
$ws.Range('E21').Value = '  +2.93%  '
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.32'
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  +18.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.53'
$ws.Range('E27').Value = '  +7.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.95'
$ws.Range('E28').Value = '  +1.38%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.977'
$ws.Range('E30').Value = '  +24.78%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0572'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.11'
$ws.Range('E32').Value = '  +2.48%  '
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.15'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('E35').Value = '  +6.70%  '
$ws.Range('E36').Value = '  +9.51%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  +2.33%  '
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '90.81'
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0625'
$ws.Range('E41').Value = '  +15.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.68'
$ws.Range('E42').Value = '  +4.05%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '50.55'
$ws.Range('E43').Value = '  +45.26%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.349.83'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.37'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.94'
$ws.Range('E46').Value = '  +1.14%  '
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  +3.19%  '
$ws.Range('D50').Value = '2.090.31'
$ws.Range('E50').Value = '  +2.71%  '
$ws.Range('E51').Value = '  +1.29%  '
